$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 3 and 4 (the original row 4 becomes row 3 after the first
# delete, so deleting row 3 twice removes both original rows 3 and 4).
$ws.Rows.Item(3).Delete()
$ws.Rows.Item(3).Delete()
